$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("gStation")

# Insert three new rows before the current row 4 (shifts old rows 4-8 down to 7-11)
$ws.Range("A4:A6").EntireRow.Insert()

# New row 4: batt.E_rated
$ws.Range("A4").Value = "batt.E_rated"
$ws.Range("B4").Value = 1.1376999999999999

# New row 5: batt.E_ex
$ws.Range("A5").Value = "batt.E_ex"
$ws.Range("B5").NumberFormat = "0.00E+00"
$ws.Range("B5").Value = 148930

# New row 6: batt.f_repl
$ws.Range("A6").Value = "batt.f_repl"
$ws.Range("B6").Value = -1

# Match the author's final cursor position recorded in the saved file
$ws.Range("D17").Select() | Out-Null
